# Convert illusion size data into dva (degrees of visual angle).
#
# The sheet holds raw illusion-size measurements (one header row of labels
# in row 1, followed by numeric data in A2:H17). Each numeric value is
# converted in place to degrees of visual angle (dva) using the standard
# small-angle visual-angle formula:
#
#     dva = DEGREES(ATAN(size / D))
#
# where D is the (fixed) viewing-distance-equivalent expressed in the same
# units as the raw size values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$D = 79200.0 / 67.0

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

# Row 1 is the header (text labels) - leave it untouched and only convert
# the numeric data rows beneath it.
$dataFirstRow = $firstRow + 1

for ($r = $dataFirstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value2
        if ($old -is [double] -or $old -is [int]) {
            $new = $excel.WorksheetFunction.Degrees($excel.WorksheetFunction.Atan($old / $D))
            $cell.Value2 = $new
        }
    }
}
